# Insert the new notes paragraphs at the top of the document (marking
# bullet headings bold, carrying forward the journal-style notes, and
# relocating the "_GoBack" bookmark into the newly inserted paragraph),
# matching the authored diff.

$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the end of the document. The
# content we are about to insert carries its own "_GoBack" bookmark, so
# remove the old one *first* (while it is still the only "_GoBack" in the
# document) to avoid any ambiguity over which same-named bookmark a later
# delete-by-name call would remove.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Create an empty anchor paragraph immediately before the existing first
# paragraph; InsertXML will replace its (empty) content with the full block
# of new paragraphs below, which keeps run- and paragraph-mark-level
# formatting (bold + bCs, sz20/szCs20) exactly as specified.
$firstParaRange = $d.Paragraphs(1).Range
$firstParaRange.InsertParagraphBefore()
$anchor = $d.Paragraphs(1).Range

$newContentXml = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Diagram of the final generative model architecture, and scientific discussions on its design [40]</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Recognisability of the best output [20]</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Realism of the best output [20]</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Unique of the outputs [20]</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>So 60 of the marks are simply for good images. I’m going to want a more robust way to produce my images, therefore.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>I may as well train on the test dataset. Nah. What I’m struggling to get my head around is identifying good Pegasus. Do we actually want those that the system cannot identify? I think it’s better to, for our training dataset, save every image individually, with an alpha value. Then I’ll have to judge them subjectively.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$anchor.InsertXML($newContentXml)
